# Consolidate text runs into a single run per paragraph for the slide
# titles (and a couple of "image" captions) across the deck. PowerPoint's
# automation layer collapses same-text assignments into a no-op, so each
# target text range is first set to an intermediate placeholder value and
# then to the final consolidated text; this forces the writer to rebuild
# the paragraph with a single run (and keeps <a:rPr/> empty, matching the
# original per-run formatting).

function Set-ConsolidatedText($shape, [string]$finalText) {
    $tr = $shape.TextFrame.TextRange
    $tr.Text = "~~tmp~~"
    $tr.Text = $finalText
}

$p = $ppt.ActivePresentation

$titles = @{
    1  = "Slide 1 (Content)"
    2  = "Slide 2 (Content)"
    3  = "Slide 3 (Content)"
    4  = "Slide 4 (Content)"
    5  = "Slide 5 (Two Content)"
    6  = "Slide 6 (Two Content Right)"
    7  = "Slide 7 (Content with Caption)"
    8  = "Slide 8 (Comparison)"
    9  = "Slide 9 (Content)"
    10 = "Slide 10 (Content)"
    11 = "Slide 11 (Content)"
    12 = "Slide 12 (Content)"
}

foreach ($slideIdx in $titles.Keys) {
    $s = $p.Slides.Item($slideIdx)
    $titleShape = $s.Shapes.Item(1)
    Set-ConsolidatedText $titleShape $titles[$slideIdx]
}

# "an image" / "An image" caption textboxes
Set-ConsolidatedText $p.Slides.Item(6).Shapes.Item(3) "an image"
Set-ConsolidatedText $p.Slides.Item(7).Shapes.Item(4) "An image"
Set-ConsolidatedText $p.Slides.Item(8).Shapes.Item(4) "An image"
